$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-19 (columns A-T), representing the updated
# TPM-based NATMI Has2-Cd44 ligand-receptor edge statistics.
$data = @(
    @("ECs", "Has2", "Cd44", "ECs", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 43.70915533333334, 131.127466, 0.0248563228721197, 0.02485632287211969, 6.132627608761778, 55.193648478856, 0.0000665202978350818, 0.00006652029783508179),
    @("ECs", "Has2", "Cd44", "FAPs", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 167.4277343333333, 502.283203, 0.09521203945945574, 0.09521203945945574, 23.49100407488312, 211.419036673948, 0.0002548057190483559, 0.0002548057190483559),
    @("ECs", "Has2", "Cd44", "Inflammatory-Mac", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 456.387034, 1369.161102, 0.259536094560534, 0.259536094560534, 64.03353493438134, 576.301814409432, 0.0006945684765177175, 0.0006945684765177175),
    @("ECs", "Has2", "Cd44", "MuSCs", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 113.6484576666667, 340.945373, 0.06462908596924523, 0.06462908596924523, 15.94548473574089, 143.509362621668, 0.0001729598569185578, 0.0001729598569185578),
    @("ECs", "Has2", "Cd44", "Neutrophils", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 750.176575, 2250.529725, 0.4266069892474147, 0.4266069892474147, 105.2537744142333, 947.2839697280999, 0.001141682304710324, 0.001141682304710324),
    @("ECs", "Has2", "Cd44", "Resolving-Mac", 1.0, 0.3333333333333333, 0.1403053333333333, 0.420916, 0.002676192217864005, 0.002676192217864005, 3.0, 1.0, 227.123347, 681.370041, 0.1291594678912306, 0.1291594678912306, 31.86661690861733, 286.799552177556, 0.0003456555628339671, 0.0003456555628339671),
    @("FAPs", "Has2", "Cd44", "ECs", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 43.70915533333334, 131.127466, 0.0248563228721197, 0.02485632287211969, 2109.191730425553, 18982.72557382997, 0.02287829476206003, 0.02287829476206003),
    @("FAPs", "Has2", "Cd44", "FAPs", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 167.4277343333333, 502.283203, 0.09521203945945574, 0.09521203945945574, 8079.249987941193, 72713.24989147072, 0.08763521116365991, 0.08763521116365991),
    @("FAPs", "Has2", "Cd44", "Inflammatory-Mac", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 456.387034, 1369.161102, 0.259536094560534, 0.259536094560534, 22023.02356669301, 198207.2121002371, 0.2388826095999059, 0.2388826095999059),
    @("FAPs", "Has2", "Cd44", "MuSCs", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 113.6484576666667, 340.945373, 0.06462908596924523, 0.06462908596924523, 5484.123068911096, 49357.10762019987, 0.05948600227853486, 0.05948600227853486),
    @("FAPs", "Has2", "Cd44", "Neutrophils", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 750.176575, 2250.529725, 0.4266069892474147, 0.4266069892474147, 36199.88115264038, 325798.9303737634, 0.3926582583341303, 0.3926582583341303),
    @("FAPs", "Has2", "Cd44", "Resolving-Mac", 3.0, 1.0, 48.25514733333333, 144.765442, 0.9204215313650299, 0.9204215313650299, 3.0, 1.0, 227.123347, 681.370041, 0.1291594678912306, 0.1291594678912306, 10959.87057232479, 98638.83515092313, 0.1188811552267389, 0.1188811552267389),
    @("MuSCs", "Has2", "Cd44", "ECs", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 43.70915533333334, 131.127466, 0.0248563228721197, 0.02485632287211969, 176.2253923257394, 1586.028530931654, 0.001911507812224586, 0.001911507812224585),
    @("MuSCs", "Has2", "Cd44", "FAPs", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 167.4277343333333, 502.283203, 0.09521203945945574, 0.09521203945945574, 675.0306187363065, 6075.275568626757, 0.00732202257674748, 0.007322022576747479),
    @("MuSCs", "Has2", "Cd44", "Inflammatory-Mac", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 456.387034, 1369.161102, 0.259536094560534, 0.259536094560534, 1840.048921231282, 16560.44029108154, 0.01995891648411038, 0.01995891648411037),
    @("MuSCs", "Has2", "Cd44", "MuSCs", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 113.6484576666667, 340.945373, 0.06462908596924523, 0.06462908596924523, 458.2047831121097, 4123.843048008987, 0.004970123833791812, 0.004970123833791811),
    @("MuSCs", "Has2", "Cd44", "Neutrophils", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 750.176575, 2250.529725, 0.4266069892474147, 0.4266069892474147, 3024.541660317475, 27220.87494285727, 0.03280704860857411, 0.0328070486085741),
    @("MuSCs", "Has2", "Cd44", "Resolving-Mac", 3.0, 1.0, 4.031773, 12.095319, 0.07690227641710612, 0.0769022764171061, 3.0, 1.0, 227.123347, 681.370041, 0.1291594678912306, 0.1291594678912306, 915.7097781042311, 8241.38800293808, 0.00993265710165776, 0.009932657101657758)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}
